# Correcting part 3 and 4: update N_Rand, Z_Score, P-value, CREAL and
# Uniqueness values for the ec_pdc_20_OUT_size3 results table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 346
$ws.Range("D2").Value = "561.9+-13.2"
$ws.Range("E2").Value = -16.3
$ws.Range("G2").Value = 97.68000000000001

# Row 3
$ws.Range("C3").Value = 799
$ws.Range("D3").Value = "893.3+-17.1"
$ws.Range("E3").Value = -5.5
$ws.Range("G3").Value = 225.58
$ws.Range("H3").Value = 13

# Row 4
$ws.Range("C4").Value = 185
$ws.Range("D4").Value = "276.0+-10.9"
$ws.Range("E4").Value = -8.359999999999999
$ws.Range("G4").Value = 52.23
$ws.Range("H4").Value = 8

# Row 5
$ws.Range("C5").Value = 1049
$ws.Range("D5").Value = "1227.1+-13.5"
$ws.Range("E5").Value = -13.19
$ws.Range("G5").Value = 296.16

# Row 6
$ws.Range("C6").Value = 354
$ws.Range("D6").Value = "222.3+-12.4"
$ws.Range("E6").Value = 10.64
$ws.Range("G6").Value = 99.94
$ws.Range("H6").Value = 12

# Row 7
$ws.Range("C7").Value = 74
$ws.Range("D7").Value = "27.6+-4.7"
$ws.Range("E7").Value = 9.81
$ws.Range("G7").Value = 20.89
$ws.Range("H7").Value = 7

# Row 8
$ws.Range("C8").Value = 455
$ws.Range("D8").Value = "621.6+-14.1"
$ws.Range("E8").Value = -11.84
$ws.Range("G8").Value = 128.46
$ws.Range("H8").Value = 11

# Row 9
$ws.Range("C9").Value = 32
$ws.Range("D9").Value = "89.1+-5.9"
$ws.Range("E9").Value = -9.699999999999999
$ws.Range("G9").Value = 9.029999999999999

# Row 10
$ws.Range("C10").Value = 12
$ws.Range("D10").Value = "20.1+-4.4"
$ws.Range("E10").Value = -1.84
$ws.Range("F10").Value = 0.9740000000000001
$ws.Range("G10").Value = 3.39
$ws.Range("H10").Value = 5

# Row 11
$ws.Range("C11").Value = 30
$ws.Range("D11").Value = "43.0+-6.1"
$ws.Range("E11").Value = -2.12
$ws.Range("F11").Value = 0.986
$ws.Range("G11").Value = 8.470000000000001
$ws.Range("H11").Value = 7

# Row 12
$ws.Range("C12").Value = 146
$ws.Range("D12").Value = "61.9+-6.6"
$ws.Range("E12").Value = 12.68
$ws.Range("G12").Value = 41.22
$ws.Range("H12").Value = 9

# Row 13
$ws.Range("C13").Value = 41
$ws.Range("D13").Value = "29.7+-4.6"
$ws.Range("E13").Value = 2.49
$ws.Range("F13").Value = 0.005
$ws.Range("G13").Value = 11.58

# Row 14
$ws.Range("C14").Value = 19
$ws.Range("D14").Value = "3.7+-1.7"
$ws.Range("E14").Value = 9.23
$ws.Range("F14").Value = 0
$ws.Range("G14").Value = 5.36
$ws.Range("H14").Value = 7
